# edit.ps1
# Applies the three changes from the commit diff:
#  1. Insert six new answer paragraphs (SOLID principles) right before the
#     empty paragraph that follows "...global overview of your project."
#  2. Insert a new GRASP-patterns answer paragraph right after the
#     "Which GRASP patterns have you used..." question paragraph.
#  3. Remove a stray Courier-New single-space run at the start of the
#     "How did you design your test suite?" paragraph.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Insert the six SOLID-principles paragraphs before the empty
#    paragraph that follows the SOLID-principles question.
# ---------------------------------------------------------------------
$anchor1 = $d.Content
$anchor1.Find.Execute("diagram representing a global overview of your project.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$questionPara1 = $anchor1.Paragraphs(1)
$targetPara1 = $questionPara1.Next()

$solidBodyXml = @'
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">To begin, we took the SOILD principles and applied them with high priority into our code. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">First, we took on the single responsibility principle. Here we took time as to go through each feature that needs to be added. What we did was create a separate class for each feature. In our code we have a </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>MeshAttributes</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> folder which contains each separate attribute that our mesh can contain. Each class in this folder tackles one feature and if there is anything that is needed that is outside the given class, we make sure we appropriately relate the two classes using the skills we learned with UML design. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">Second, we look at the Open Closes Principle. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">Third, we have </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>Liskov</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> Substitution Principle</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">. Here we can see this used in our interface called </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>BoundedShapes</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">. Each shape that extends this interface can easily be switched out for our </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>parent. For bounded shape, we have our bounded method and scale method, where bounded is the coordinates of or shape and scale is the scaled factor. We then have square and circle that extend our interface. Each one can be switched and replaced</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> to follow the principle as our shape is abstract and can contain any shape.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>Fourth</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>, we have interface segregation principle.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> We made sure that any functions not needed by our user were not easily accessible. For example, we have either an irregular, square or circular shape. This meant that if you were to access the interface that contains this class, you would not be able to access any of the </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>functions that</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> do not concern the user.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:lang w:val="en-CA"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">Finally, we have the Dependency Principle. Here we made our abstractions as our key priority. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>We realized that our code is depends too much on the visualizer class and a lot of the things in our MVP was hard coded. We decided to make sure that everything can be changed to so that it is abstract instead of being hard coded with if statements. This is seen through our lakes and rivers in our generated mesh. We made it so if the points in our bounded shape were inside a certain area, they would be lakes. We quickly saw this didn’t rely on abstraction and quickly changed this in our code.</w:t>
      </w:r>
    </w:p>
'@

$solidPkgXml = New-PkgXml($solidBodyXml)
$insertionPoint1 = $targetPara1.Range
$insertionPoint1.Collapse(1)
$insertionPoint1.InsertXML($solidPkgXml)

# ---------------------------------------------------------------------
# 2) Insert the GRASP-patterns answer paragraph right after the
#    "Which GRASP patterns..." question paragraph.
# ---------------------------------------------------------------------
$anchor2 = $d.Content
$anchor2.Find.Execute("Which GRASP patterns have you used when attributing responsibilities?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$graspQuestionPara = $anchor2.Paragraphs(1)
$targetPara2 = $graspQuestionPara.Next()

$graspBodyXml = @'
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
        <w:rPr>
          <w:lang w:val="en-CA"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t xml:space="preserve">The different GRASP patterns we used were High Cohesion, Low Coupling, and Polymorphism. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-CA"/>
        </w:rPr>
        <w:t xml:space="preserve">For high cohesion, we designed each class in our program to have its own purpose. We made it so that every single has one responsibility and has a clear focus and purpose. This makes it so changes for testing, changing, or adding features is easy and reliable.  </w:t>
      </w:r>
    </w:p>
'@

$graspPkgXml = New-PkgXml($graspBodyXml)
$insertionPoint2 = $targetPara2.Range
$insertionPoint2.Collapse(1)
$insertionPoint2.InsertXML($graspPkgXml)

# ---------------------------------------------------------------------
# 3) Remove the stray Courier-New single-space run at the start of the
#    "How did you design your test suite?" paragraph.
# ---------------------------------------------------------------------
$anchor3 = $d.Content
$anchor3.Find.Execute("How did you design your test suite?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$testSuitePara = $anchor3.Paragraphs(1)
$paraRange = $testSuitePara.Range
$paraStart = $paraRange.Start
$spaceRange = $d.Range($paraStart, $paraStart + 1)
Write-Output ("Space-run candidate text=[" + $spaceRange.Text + "] font=" + $spaceRange.Font.Name)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

Write-Output "Done."
